$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates scraped from the coinranking.com refresh (GitHub Actions bot).
# Column D ("Price") holds numeric-looking text (e.g. "26.826.52", "1.000") that must
# stay literal text, not be auto-coerced to a number by Excel. We write it via a
# text formula + Copy/PasteSpecial(xlPasteValues) round-trip, which is how Excel
# itself converts a formula result into a static value while keeping its string type.
$xlPasteValues = -4163

function Set-PriceText($row, $text) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteValues) | Out-Null
}

# Row 2
Set-PriceText 2 '26.826.52'
$ws.Cells.Item(2, 5).Value = '  -1.42%  '

# Row 3
Set-PriceText 3 '1.872.84'
$ws.Cells.Item(3, 5).Value = '  -1.69%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.24%  '

# Row 5
Set-PriceText 5 '300.88'
$ws.Cells.Item(5, 5).Value = '  -2.18%  '

# Row 6
Set-PriceText 6 '1.000'
$ws.Cells.Item(6, 5).Value = '  -0.20%  '

# Row 7
Set-PriceText 7 '0.5326'
$ws.Cells.Item(7, 5).Value = '  +1.33%  '

# Row 8
Set-PriceText 8 '0.3758'
$ws.Cells.Item(8, 5).Value = '  -1.43%  '

# Row 9
Set-PriceText 9 '0.07179'
$ws.Cells.Item(9, 5).Value = '  -1.71%  '

# Row 10
Set-PriceText 10 '21.63'
$ws.Cells.Item(10, 5).Value = '  -0.39%  '

# Row 11
Set-PriceText 11 '0.8868'
$ws.Cells.Item(11, 5).Value = '  -2.01%  '

# Row 12
Set-PriceText 12 '0.08167'
$ws.Cells.Item(12, 5).Value = '  +1.42%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-PriceText 13 '1.875.06'
$ws.Cells.Item(13, 5).Value = '  +4.27%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Litecoin'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-PriceText 14 '93.39'
$ws.Cells.Item(14, 5).Value = '  -2.63%  '

# Row 15
Set-PriceText 15 '5.282'
$ws.Cells.Item(15, 5).Value = '  -1.54%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.19%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.28%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  -1.63%  '

# Row 19
Set-PriceText 19 '1.000'
$ws.Cells.Item(19, 5).Value = '  -0.17%  '

# Row 20
Set-PriceText 20 '26.813.73'
$ws.Cells.Item(20, 5).Value = '  -1.60%  '

# Row 21
Set-PriceText 21 '4.980'
$ws.Cells.Item(21, 5).Value = '  -2.82%  '

# Row 22
Set-PriceText 22 '10.67'
$ws.Cells.Item(22, 5).Value = '  -1.40%  '

# Row 23
Set-PriceText 23 '6.385'
$ws.Cells.Item(23, 5).Value = '  -1.45%  '

# Row 24
Set-PriceText 24 '146.36'
$ws.Cells.Item(24, 5).Value = '  -2.02%  '

# Row 25
Set-PriceText 25 '2.269'
$ws.Cells.Item(25, 5).Value = '  -3.54%  '

# Row 26
Set-PriceText 26 '1.732'
$ws.Cells.Item(26, 5).Value = '  -0.61%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -1.30%  '

# Row 28
Set-PriceText 28 '113.87'
$ws.Cells.Item(28, 5).Value = '  -2.62%  '

# Row 29
Set-PriceText 29 '4.724'
$ws.Cells.Item(29, 5).Value = '  -2.50%  '

# Row 30
Set-PriceText 30 '4.617'
$ws.Cells.Item(30, 5).Value = '  -5.86%  '

# Row 31
Set-PriceText 31 '0.09132'
$ws.Cells.Item(31, 5).Value = '  -1.19%  '

# Row 32
Set-PriceText 32 '0.8071'
$ws.Cells.Item(32, 5).Value = '  +0.49%  '

# Row 33
Set-PriceText 33 '0.04972'
$ws.Cells.Item(33, 5).Value = '  -2.19%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -4.50%  '

# Row 35
Set-PriceText 35 '2.960'
$ws.Cells.Item(35, 5).Value = '  -0.31%  '

# Row 36
Set-PriceText 36 '0.6016'
$ws.Cells.Item(36, 5).Value = '  +5.04%  '

# Row 37
Set-PriceText 37 '3.191'
$ws.Cells.Item(37, 5).Value = '  -5.66%  '

# Row 38
Set-PriceText 38 '2.589'
$ws.Cells.Item(38, 5).Value = '  -2.95%  '

# Row 39
Set-PriceText 39 '0.01949'
$ws.Cells.Item(39, 5).Value = '  -2.23%  '

# Row 40
Set-PriceText 40 '1.070'
$ws.Cells.Item(40, 5).Value = '  -1.66%  '

# Row 41
Set-PriceText 41 '8.887'
$ws.Cells.Item(41, 5).Value = '  -1.21%  '

# Row 42
Set-PriceText 42 '6.564'
$ws.Cells.Item(42, 5).Value = '  -0.63%  '

# Row 43
Set-PriceText 43 '0.5134'
$ws.Cells.Item(43, 5).Value = '  +4.60%  '

# Row 44
Set-PriceText 44 '114.37'
$ws.Cells.Item(44, 5).Value = '  -2.11%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.89%  '

# Row 46
Set-PriceText 46 '0.9999'
$ws.Cells.Item(46, 5).Value = '  -0.25%  '

# Row 47
Set-PriceText 47 '9.933'
$ws.Cells.Item(47, 5).Value = '  -2.35%  '

# Row 48
Set-PriceText 48 '1.634'
$ws.Cells.Item(48, 5).Value = '  -0.34%  '

# Row 49
Set-PriceText 49 '37.50'
$ws.Cells.Item(49, 5).Value = '  -2.86%  '

# Row 50
Set-PriceText 50 '0.06047'
$ws.Cells.Item(50, 5).Value = '  +1.43%  '

# Row 51
Set-PriceText 51 '61.99'
$ws.Cells.Item(51, 5).Value = '  -3.90%  '
